# Update the risk-factor multipliers on "prevalence2018" and propagate the
# knock-on recalculation of "incidence2018_plus" (which derives its values
# from prevalence2018 via a formula).

$wb = $excel.ActiveWorkbook

$wsPrev = $wb.Worksheets.Item("prevalence2018")
$wsInc  = $wb.Worksheets.Item("incidence2018_plus")

# --- prevalence2018: new plateau values for column C, rows 27-122 ---------
$wsPrev.Range("C27:C36").Value  = 0.211
$wsPrev.Range("C37:C46").Value  = 0.314
$wsPrev.Range("C47:C56").Value  = 0.429
$wsPrev.Range("C57:C122").Value = 0.593

# --- incidence2018_plus: formulas divide by 2 instead of by 8 -------------
for ($r = 27; $r -le 122; $r++) {
    $wsInc.Range("C$r").Formula = "=prevalence2018!C$r/2"
}

# --- sheet view / selection changes ---------------------------------------
$wsPrev.Activate()
$excel.ActiveWindow.ScrollRow = 21
$wsPrev.Range("E27:E31").Select()

$wsInc.Range("H29").Select()
